$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 359-360. This shifts the previous rows 359..377
# down to become rows 361..379 (unchanged content), matching the diff's
# dimension growth from A1:R377 to A1:R379.
$ws.Range("A359:R360").Insert()

# New row 359: Cebolla, 1a (cosecha), week of 44610
$ws.Cells.Item(359, 1).Value = 11
$ws.Cells.Item(359, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(359, 3).Value = "Bíobío"
$ws.Cells.Item(359, 4).Value = 44610
$ws.Cells.Item(359, 5).Value = 8
$ws.Cells.Item(359, 6).Value = 100112004
$ws.Cells.Item(359, 7).Value = "Cebolla"
$ws.Cells.Item(359, 8).Value = "Sin especificar"
$ws.Cells.Item(359, 9).Value = "1a (cosecha)"
$ws.Cells.Item(359, 10).Value = 1000
$ws.Cells.Item(359, 11).Value = 4500
$ws.Cells.Item(359, 12).Value = 5000
$ws.Cells.Item(359, 13).Value = 4750
$ws.Cells.Item(359, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(359, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(359, 16).Value = 264
$ws.Cells.Item(359, 17).Value = 18
$ws.Cells.Item(359, 18).Value = "Hortaliza"

# New row 360: Cebolla, 2a (cosecha), week of 44610
$ws.Cells.Item(360, 1).Value = 11
$ws.Cells.Item(360, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(360, 3).Value = "Bíobío"
$ws.Cells.Item(360, 4).Value = 44610
$ws.Cells.Item(360, 5).Value = 8
$ws.Cells.Item(360, 6).Value = 100112004
$ws.Cells.Item(360, 7).Value = "Cebolla"
$ws.Cells.Item(360, 8).Value = "Sin especificar"
$ws.Cells.Item(360, 9).Value = "2a (cosecha)"
$ws.Cells.Item(360, 10).Value = 500
$ws.Cells.Item(360, 11).Value = 4000
$ws.Cells.Item(360, 12).Value = 4000
$ws.Cells.Item(360, 13).Value = 4000
$ws.Cells.Item(360, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(360, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(360, 16).Value = 222
$ws.Cells.Item(360, 17).Value = 18
$ws.Cells.Item(360, 18).Value = "Hortaliza"
